$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = "yahuza"
$ws.Range("C2").Value = "umar"
$ws.Range("D2").Value = 1234567
$ws.Range("E2").Value = "male"
$ws.Range("F2").Value = 1000
$ws.Range("G2").Value = "geology"

# Row 3
$ws.Range("B3").Value = "aisha"
$ws.Range("C3").Value = "hassan"
$ws.Range("D3").Value = 3049343
$ws.Range("E3").Value = "female"
$ws.Range("F3").Value = 334
$ws.Range("G3").Value = "fam"

# First Name column filled last
$ws.Range("A2").Value = "attama"
$ws.Range("A3").Value = "oga"

$ws.Range("A3").Select()
